# Task & work time documentation.xlsx - apply commit changes
# "changed readme.md, tried to fix ref issue / ref issue was NOT resolved added time"
#
# Content changes applied here:
#  - New work-log row (J3:L3) on Tabelle1: task text, duration, date
#  - Selection cursor moved to J10 (reflects where the user ended up working)
#  - Column widths for H, J, K, L nudged to fit the new content
#  - (Cosmetic note: the source diff also touches a couple of pure
#    environment/session artifacts - the absPath/revisionPtr GUIDs in
#    workbook.xml and the built-in "Normal" cellStyle display name - which
#    are written by Excel's own save pipeline from the local machine/locale
#    and are not exposed anywhere on the Workbook/Worksheet/Range object
#    model, so they are not reachable from this script.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data: Task / Time(h) / Date, mirroring the B:D and F:H blocks ---
$ws.Range("J3").Value = "fixing ref / dependency issues of poppinger"
$ws.Range("K3").Value = "2h"
$ws.Range("L3").Value = 45678
$ws.Range("L3").NumberFormat = "d-mmm"

# --- Column widths: H/J/K/L need to widen to comfortably show the new entries ---
$ws.Columns.Item(8).ColumnWidth = 15.45
$ws.Columns.Item(10).ColumnWidth = 18.3
$ws.Columns.Item(11).ColumnWidth = 14.17
$ws.Columns.Item(12).ColumnWidth = 15.17

# --- Selection left on J10 ---
$ws.Range("J10").Select() | Out-Null

Write-Output "edit applied"
